$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F50: add centered style (matches xf index 1 in the original style table) ---
$ws.Range("F50").HorizontalAlignment = -4108

# --- F51: was using the lone "numFmt 164, no alignment" style (xf 5); that
#     style is being retired, so F51 moves onto the "numFmt 164 + centered"
#     style (xf 3) already used by the rest of the F43:F49 column. ---
$ws.Range("F51").NumberFormat = "0.000"
$ws.Range("F51").HorizontalAlignment = -4108

# --- New row 52: STDEV summary row under the last data block ---
$ws.Range("A52").Value = "STDEV"
$ws.Range("B52").Formula = "=STDEV(B43:B50)"
$ws.Range("C52").Formula = "=STDEV(C43:C50)"
$ws.Range("D52").Formula = "=STDEV(D43:D50)"
$ws.Range("E52").Formula = "=STDEV(E43:E50)"
$ws.Range("F52").Formula = "=STDEV(F43:F51)"

$ws.Range("A52:D52").NumberFormat = "0.000"
$ws.Range("A52:D52").HorizontalAlignment = -4108
$ws.Range("F52").NumberFormat = "0.000"
$ws.Range("F52").HorizontalAlignment = -4108

# --- View state: scroll so row 12 is at top, select the new summary block ---
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A51:F52").Select()
